$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.42%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.99%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.587"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.16%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08101"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.23%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.981"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.52%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D7").Value = "'2.574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-4.88%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.42%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1175"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-3.01%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1862"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.35%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D11").Value = "'10.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'19.79%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09854"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'2.68%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04655"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'12.43%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1068"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.03%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001283"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.99%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04220"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-3.52%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006007"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.75%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.372"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-5.76%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'4.320"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.64%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3475"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.68%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1409"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'4.61%"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2508"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.58%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "'0.001250"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.15%"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "'0.004325"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.78%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001192"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-3.58%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-0.57%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02661"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'0.22%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05568"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'2.26%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007561"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-2.08%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1407"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.09%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.008083"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-16.97%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002019"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-4.83%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008413"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-14.77%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00007229"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.59%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.17%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.004733"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'33.02%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002273"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-0.18%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.17%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002003"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.17%"
$ws.Range("E50").Style = "Normal"
